$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-150 down to 42-151
$ws.Rows(41).Insert()

# Fill new row 41 with data
$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C41").Value = "Los Lagos"
$ws.Range("D41").Value = 44414
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = 100112045
$ws.Range("G41").Value = "Zapallo"
$ws.Range("H41").Value = "Paine"
$ws.Range("I41").Value = "1a (guarda)"
$ws.Range("J41").Value = 1100
$ws.Range("K41").Value = 450
$ws.Range("L41").Value = 450
$ws.Range("M41").Value = 450
$ws.Range("N41").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O41").Value = "Región de O'Higgins"
$ws.Range("P41").Value = 450
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
